$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the existing row 228, shifting the
# remaining rows (old 228-237) down to 230-239.
$ws.Rows.Item(228).Insert()
$ws.Rows.Item(229).Insert()

# New row 228
$ws.Cells.Item(228, 1).Value = 9
$ws.Cells.Item(228, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(228, 3).Value = "Metropolitana"
$ws.Cells.Item(228, 4).Value = 44610
$ws.Cells.Item(228, 5).Value = 13
$ws.Cells.Item(228, 6).Value = 100112030
$ws.Cells.Item(228, 7).Value = "Poroto granado"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 61
$ws.Cells.Item(228, 11).Value = 24000
$ws.Cells.Item(228, 12).Value = 26000
$ws.Cells.Item(228, 13).Value = 25016
$ws.Cells.Item(228, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(228, 15).Value = "Región Metropolitana"
$ws.Cells.Item(228, 16).Value = 1001
$ws.Cells.Item(228, 17).Value = 25
$ws.Cells.Item(228, 18).Value = "Hortaliza"

# New row 229
$ws.Cells.Item(229, 1).Value = 9
$ws.Cells.Item(229, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(229, 3).Value = "Metropolitana"
$ws.Cells.Item(229, 4).Value = 44610
$ws.Cells.Item(229, 5).Value = 13
$ws.Cells.Item(229, 6).Value = 100112030
$ws.Cells.Item(229, 7).Value = "Poroto granado"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 52
$ws.Cells.Item(229, 11).Value = 24000
$ws.Cells.Item(229, 12).Value = 26000
$ws.Cells.Item(229, 13).Value = 25000
$ws.Cells.Item(229, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(229, 15).Value = "Región del Maule"
$ws.Cells.Item(229, 16).Value = 1000
$ws.Cells.Item(229, 17).Value = 25
$ws.Cells.Item(229, 18).Value = "Hortaliza"
